$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value2 = 715
$ws.Cells.Item(3, 6).Value2 = 13796
$ws.Cells.Item(4, 6).Value2 = 13798
$ws.Cells.Item(5, 6).Value2 = 13882
$ws.Cells.Item(7, 6).Value2 = 1355
$ws.Cells.Item(8, 6).Value2 = 1371
$ws.Cells.Item(9, 6).Value2 = 5769
$ws.Cells.Item(10, 6).Value2 = 963
$ws.Cells.Item(15, 6).Value2 = 1507
$ws.Cells.Item(16, 6).Value2 = 415
$ws.Cells.Item(18, 6).Value2 = 1157
$ws.Cells.Item(19, 6).Value2 = 1763
$ws.Cells.Item(20, 6).Value2 = 903
$ws.Cells.Item(21, 6).Value2 = 28
$ws.Cells.Item(22, 6).Value2 = 2235
$ws.Cells.Item(24, 6).Value2 = 781
$ws.Cells.Item(25, 6).Value2 = 3222
$ws.Cells.Item(27, 6).Value2 = 296
$ws.Cells.Item(28, 6).Value2 = 2279
$ws.Cells.Item(29, 6).Value2 = 68
$ws.Cells.Item(31, 6).Value2 = 1339
$ws.Cells.Item(32, 6).Value2 = 1749
$ws.Cells.Item(33, 6).Value2 = 1054
$ws.Cells.Item(34, 6).Value2 = 1303
$ws.Cells.Item(35, 6).Value2 = 89
$ws.Cells.Item(36, 6).Value2 = 129
$ws.Cells.Item(37, 6).Value2 = 4580
$ws.Cells.Item(38, 6).Value2 = 4688
$ws.Cells.Item(40, 6).Value2 = 150
$ws.Cells.Item(41, 6).Value2 = 655
$ws.Cells.Item(42, 6).Value2 = 667
$ws.Cells.Item(43, 6).Value2 = 3245
$ws.Cells.Item(47, 6).Value2 = 78
$ws.Cells.Item(49, 6).Value2 = 4388
$ws.Cells.Item(50, 6).Value2 = 267

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value2 = 107

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value2 = 7359
$ws.Cells.Item(3, 6).Value2 = 200
$ws.Cells.Item(4, 6).Value2 = 627

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value2 = 7359
$ws.Cells.Item(3, 6).Value2 = 715
$ws.Cells.Item(4, 6).Value2 = 200
$ws.Cells.Item(5, 6).Value2 = 627
$ws.Cells.Item(7, 6).Value2 = 13800
$ws.Cells.Item(8, 6).Value2 = 13883
$ws.Cells.Item(9, 6).Value2 = 1355
$ws.Cells.Item(10, 6).Value2 = 1371
$ws.Cells.Item(11, 6).Value2 = 5769
$ws.Cells.Item(17, 6).Value2 = 415
$ws.Cells.Item(18, 6).Value2 = 1157
$ws.Cells.Item(19, 6).Value2 = 1763
$ws.Cells.Item(21, 6).Value2 = 781
$ws.Cells.Item(22, 6).Value2 = 3222
$ws.Cells.Item(23, 6).Value2 = 296
$ws.Cells.Item(24, 6).Value2 = 68
$ws.Cells.Item(27, 6).Value2 = 1749
$ws.Cells.Item(33, 6).Value2 = 1054
$ws.Cells.Item(34, 6).Value2 = 1303
$ws.Cells.Item(36, 6).Value2 = 4580
$ws.Cells.Item(37, 6).Value2 = 4688
$ws.Cells.Item(39, 6).Value2 = 150
$ws.Cells.Item(40, 6).Value2 = 3245
$ws.Cells.Item(44, 6).Value2 = 78
$ws.Cells.Item(46, 6).Value2 = 4388
$ws.Cells.Item(47, 6).Value2 = 267
